$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42632.878750000003
$ws.Range("B3").Value = -1
$ws.Range("C3").Value = "Neutral"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 17236
$ws.Range("F3").Value = 891
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 41
$ws.Range("I3").Value = 91
$ws.Range("J3").Value = 8
$ws.Range("K3").Value = 18436
$ws.Range("L3").Value = 202
$ws.Range("M3").Value = 142
$ws.Range("N3").Value = 21
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = "Named"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = -31.57
$ws.Range("S3").Value = -0.0872
$ws.Range("S3").NumberFormat = $ws.Range("S2").NumberFormat
$ws.Range("T3").Value = -0.74
$ws.Range("U3").Value = 6.75
$ws.Range("V3").Value = 1.88
$ws.Range("W3").Value = 0
